# Update the "dSF" (delta-S-final) column (F) values for several rows,
# as part of a data repull / push-all-data / mean-calculation refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = -1
$ws.Range("F6").Value = -8
$ws.Range("F7").Value = 11
$ws.Range("F8").Value = -7
$ws.Range("F9").Value = 0
$ws.Range("F12").Value = -11
$ws.Range("F26").Value = -2
$ws.Range("F27").Value = 4
